# Update the title text on slide 1 from "hello" to "Hello world".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$titleShape = $s.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Hello world"
